$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 454.14285
$ws.Range("I2").Value = 379.83334
$ws.Range("J2").Value = 900.0
$ws.Range("K2").Value = 379.83334
$ws.Range("L2").Value = 900.0
$ws.Range("M2").Value = -266.83334
$ws.Range("N2").Value = -1126.0
$ws.Range("H4").Value = 8408.667
$ws.Range("I4").Value = 5590.4
$ws.Range("J4").Value = 22500.0
$ws.Range("K4").Value = 5590.4
$ws.Range("L4").Value = 22500.0
$ws.Range("M4").Value = -5476.4
$ws.Range("N4").Value = -22728.0
$ws.Range("H9").Value = 171.82353
$ws.Range("I9").Value = 132.33333
$ws.Range("J9").Value = 266.6
$ws.Range("K9").Value = 132.33333
$ws.Range("L9").Value = 266.6
$ws.Range("M9").Value = 36.66667000000001
$ws.Range("N9").Value = -604.6
$ws.Range("H32").Value = 8849.8
$ws.Range("I32").Value = 2250.5
$ws.Range("K32").Value = 2250.5
$ws.Range("M32").Value = -1924.5
$ws.Range("H38").Value = 8868.385
$ws.Range("J38").Value = 50009.0
$ws.Range("L38").Value = 150027.0
$ws.Range("N38").Value = -150771.0
$ws.Range("H40").Value = 6481.25
$ws.Range("I40").Value = 4949.222
$ws.Range("K40").Value = 4949.222
$ws.Range("M40").Value = -4774.222
$ws.Range("H43").Value = 3998.5
$ws.Range("J43").Value = 3998.5
$ws.Range("L43").Value = 3998.5
$ws.Range("N43").Value = -4136.5
$ws.Range("H64").Value = 30003332.0
$ws.Range("J64").Value = 4999.5
$ws.Range("L64").Value = 4999.5
$ws.Range("N64").Value = -5495.5
$ws.Range("H67").Value = 30003332.0
$ws.Range("J67").Value = 4999.5
$ws.Range("L67").Value = 4999.5
$ws.Range("N67").Value = -6715.5
$ws.Range("H69").Value = 5257.25
$ws.Range("J69").Value = 0.0
$ws.Range("L69").Value = 0.0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 5257.25
$ws.Range("J72").Value = 0.0
$ws.Range("L72").Value = 0.0
$ws.Range("N72").ClearContents()
$ws.Range("H110").Value = 84997.5
$ws.Range("J110").Value = 99995.0
$ws.Range("L110").Value = 99995.0
$ws.Range("N110").Value = -108175.0
$ws.Range("H113").Value = 2949.5
$ws.Range("J113").Value = 2900.0
$ws.Range("L113").Value = 2900.0
$ws.Range("N113").Value = -9408.0
$ws.Range("H138").Value = 3499.5
$ws.Range("J138").Value = 4874.25
$ws.Range("L138").Value = 14622.75
$ws.Range("N138").Value = -24902.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 750.0
$ws.Range("I5").Value = 500.0
$ws.Range("K5").Value = 500.0
$ws.Range("M5").Value = -388.0
$ws.Range("H28").Value = 28700.0
$ws.Range("I28").Value = 28700.0
$ws.Range("K28").Value = 28700.0
$ws.Range("M28").Value = -28508.0
$ws.Range("H45").Value = 2475.0
$ws.Range("I45").Value = 2475.0
$ws.Range("K45").Value = 2475.0
$ws.Range("M45").Value = -2098.0
$ws.Range("H88").Value = 1966.6666
$ws.Range("I88").Value = 1160.0
$ws.Range("K88").Value = 1160.0
$ws.Range("M88").Value = -754.0
$ws.Range("H91").Value = 1966.6666
$ws.Range("I91").Value = 1160.0
$ws.Range("K91").Value = 1160.0
$ws.Range("M91").Value = 244.0
$ws.Range("H97").Value = 501.5
$ws.Range("I97").Value = 501.5
$ws.Range("K97").Value = 501.5
$ws.Range("M97").Value = -5.5
$ws.Range("H99").Value = 28700.0
$ws.Range("I99").Value = 28700.0
$ws.Range("K99").Value = 28700.0
$ws.Range("M99").Value = -25705.0
$ws.Range("H114").Value = 0.0
$ws.Range("J114").Value = 0.0
$ws.Range("L114").Value = 0.0
$ws.Range("N114").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 750.0
$ws.Range("I4").Value = 500.0
$ws.Range("K4").Value = 500.0
$ws.Range("M4").Value = -385.0
$ws.Range("H86").Value = 36534.5
$ws.Range("I86").Value = 4200.0
$ws.Range("K86").Value = 4200.0
$ws.Range("M86").Value = -3077.0
$ws.Range("H89").Value = 36534.5
$ws.Range("I89").Value = 4200.0
$ws.Range("K89").Value = 21000.0
$ws.Range("M89").Value = -15384.0
$ws.Range("H94").Value = 3552.125
$ws.Range("I94").Value = 3151.1667
$ws.Range("K94").Value = 3151.1667
$ws.Range("M94").Value = -2700.1667
$ws.Range("H99").Value = 2900.0
$ws.Range("I99").Value = 2500.0
$ws.Range("J99").Value = 3166.6667
$ws.Range("K99").Value = 2500.0
$ws.Range("L99").Value = 3166.6667
$ws.Range("M99").Value = -1002.0
$ws.Range("N99").Value = -6162.6667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 461910.47
$ws.Range("I4").Value = 417069.75
$ws.Range("K4").Value = 1251209.25
$ws.Range("M4").Value = -1251097.25
$ws.Range("H34").Value = 1643.75
$ws.Range("I34").Value = 324.0
$ws.Range("K34").Value = 972.0
$ws.Range("M34").Value = -888.0
$ws.Range("H39").Value = 3513.7144
$ws.Range("I39").Value = 2197.6667
$ws.Range("J39").Value = 4500.75
$ws.Range("K39").Value = 6593.000100000001
$ws.Range("L39").Value = 13502.25
$ws.Range("M39").Value = -6299.000100000001
$ws.Range("N39").Value = -14090.25
$ws.Range("H55").Value = 2450.4
$ws.Range("J55").Value = 3876.0
$ws.Range("L55").Value = 11628.0
$ws.Range("N55").Value = -11982.0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3421.4285
$ws.Range("I97").Value = 3408.3333
$ws.Range("J97").Value = 3500.0
$ws.Range("K97").Value = 3408.3333
$ws.Range("L97").Value = 3500.0
$ws.Range("M97").Value = -2912.3333
$ws.Range("N97").Value = -4492.0
$ws.Range("H98").Value = 39999.0
$ws.Range("J98").Value = 39999.0
$ws.Range("L98").Value = 39999.0
$ws.Range("N98").Value = -45989.0
$ws.Range("H132").Value = 6312.222
$ws.Range("I132").Value = 5687.143
$ws.Range("K132").Value = 17061.429
$ws.Range("M132").Value = -14531.429

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 5682.75
$ws.Range("I21").Value = 6303.0
$ws.Range("J21").Value = 5062.5
$ws.Range("K21").Value = 6303.0
$ws.Range("L21").Value = 5062.5
$ws.Range("M21").Value = -6129.0
$ws.Range("N21").Value = -5410.5
$ws.Range("H55").Value = 2109.0588
$ws.Range("I55").Value = 1543.25
$ws.Range("J55").Value = 2612.0
$ws.Range("K55").Value = 1543.25
$ws.Range("L55").Value = 2612.0
$ws.Range("M55").Value = -1370.25
$ws.Range("N55").Value = -2958.0
$ws.Range("H93").Value = 12048.75
$ws.Range("I93").Value = 15065.0
$ws.Range("K93").Value = 15065.0
$ws.Range("M93").Value = -13817.0
$ws.Range("H100").Value = 4125.5
$ws.Range("I100").Value = 3501.0
$ws.Range("K100").Value = 3501.0
$ws.Range("M100").Value = -2960.0
$ws.Range("H132").Value = 4391.9165
$ws.Range("J132").Value = 4244.778
$ws.Range("L132").Value = 12734.334
$ws.Range("N132").Value = -17794.334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 41056.832
$ws.Range("I2").Value = 60166.25
$ws.Range("J2").Value = 2838.0
$ws.Range("K2").Value = 60166.25
$ws.Range("L2").Value = 2838.0
$ws.Range("M2").Value = -60054.25
$ws.Range("N2").Value = -3062.0
